$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells in row 3 (P3:R3), mirroring the existing header style (O3) ---
$ws.Range("P3").Value = "DOCTOR’S NAME"
$ws.Range("Q3").Value = "MEDICAL LICENSE"
$ws.Range("R3").Value = "INSTITUTE / HOSPITAL"

$ws.Range("P3:R3").HorizontalAlignment = $ws.Range("O3").HorizontalAlignment
$ws.Range("P3:R3").VerticalAlignment = $ws.Range("O3").VerticalAlignment

# --- Column width tweaks (ColumnWidth is in characters; stored xlsx width = ColumnWidth + 5/6) ---
$ws.Columns.Item(1).ColumnWidth = 3.3366666666666664    # -> 4.17
$ws.Columns.Item(4).ColumnWidth = 5.396666666666667     # -> 6.23
$ws.Columns.Item(12).ColumnWidth = 8.336666666666666    # -> 9.17 (new column)
$ws.Columns.Item(13).ColumnWidth = 8.196666666666665    # -> 9.03
$ws.Columns.Item(14).ColumnWidth = 12.796666666666667   # -> 13.63
$ws.Columns.Item(16).ColumnWidth = 16.256666666666668   # -> 17.09 (new column)
$ws.Columns.Item(17).ColumnWidth = 17.366666666666667   # -> 18.2  (new column)
$ws.Columns.Item(18).ColumnWidth = 21.526666666666667   # -> 22.36 (new column)

# --- Selection moves to the newly added last cell ---
[void]$ws.Range("R3").Select()
